$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_4_8_0"
$ws.Range("B2").Value = 0.5004926934021611
$ws.Range("C2").Value = 0.824446744700852
$ws.Range("D2").Value = 0.8378541844271405
$ws.Range("E2").Value = 0.8315618278201017
$ws.Range("F2").Value = 0.5528073906898499
$ws.Range("G2").Value = 0.2289877682924271
$ws.Range("H2").Value = 0.1838966757059097
$ws.Range("I2").Value = 0.207768440246582

$ws.Range("A3").Value = "model_4_8_1"
$ws.Range("B3").Value = 0.5370917275399447
$ws.Range("C3").Value = 0.8216897894654159
$ws.Range("D3").Value = 0.7514891257035685
$ws.Range("E3").Value = 0.7926496294314354
$ws.Range("F3").Value = 0.5123030543327332
$ws.Range("G3").Value = 0.2325838655233383
$ws.Range("H3").Value = 0.2818470597267151
$ws.Range("I3").Value = 0.2557666301727295

$ws.Range("A4").Value = "model_4_8_21"
$ws.Range("B4").Value = 0.6216843769568885
$ws.Range("C4").Value = -0.1538938502984271
$ws.Range("D4").Value = 0.03973894339921935
$ws.Range("E4").Value = -0.06147598979003166
$ws.Range("F4").Value = 0.4186838567256927
$ws.Range("G4").Value = 1.505113363265991
$ws.Range("H4").Value = 1.08907413482666
$ws.Range("I4").Value = 1.309330224990845

$ws.Range("A5").Value = "model_4_8_22"
$ws.Range("B5").Value = 0.6246491648554175
$ws.Range("C5").Value = -0.1443440280570365
$ws.Range("D5").Value = 0.04670928445188971
$ws.Range("E5").Value = -0.05311372245393797
$ws.Range("F5").Value = 0.415402740240097
$ws.Range("G5").Value = 1.492656946182251
$ws.Range("H5").Value = 1.0811687707901
$ws.Range("I5").Value = 1.299015522003174

$ws.Range("A6").Value = "model_4_8_23"
$ws.Range("B6").Value = 0.6248814228733295
$ws.Range("C6").Value = -0.1452157333405322
$ws.Range("D6").Value = 0.04891009324202511
$ws.Range("E6").Value = -0.05264953024013574
$ws.Range("F6").Value = 0.4151457250118256
$ws.Range("G6").Value = 1.493793964385986
$ws.Range("H6").Value = 1.078672647476196
$ws.Range("I6").Value = 1.298442840576172

$ws.Range("A7").Value = "model_4_8_24"
$ws.Range("B7").Value = 0.6250513292857243
$ws.Range("C7").Value = -0.1451467909541622
$ws.Range("D7").Value = 0.04960988236461727
$ws.Range("E7").Value = -0.05230814847692677
$ws.Range("F7").Value = 0.4149576723575592
$ws.Range("G7").Value = 1.493704080581665
$ws.Range("H7").Value = 1.077878952026367
$ws.Range("I7").Value = 1.298021674156189

$ws.Range("A8").Value = "model_4_8_2"
$ws.Range("B8").Value = 0.6263881645385998
$ws.Range("C8").Value = 0.4501810612946222
$ws.Range("D8").Value = 0.5509275005255272
$ws.Range("E8").Value = 0.4978874039461524
$ws.Range("F8").Value = 0.4134782254695892
$ws.Range("G8").Value = 0.7171716094017029
$ws.Range("H8").Value = 0.5093128085136414
$ws.Range("I8").Value = 0.6193556785583496

$ws.Range("A9").Value = "model_4_8_15"
$ws.Range("B9").Value = 0.6269065238847529
$ws.Range("C9").Value = -0.1320177567459038
$ws.Range("D9").Value = 0.05106931095122302
$ws.Range("E9").Value = -0.04432676264263846
$ws.Range("F9").Value = 0.4129045307636261
$ws.Range("G9").Value = 1.476578831672668
$ws.Range("H9").Value = 1.076223850250244
$ws.Range("I9").Value = 1.288176536560059

$ws.Range("A10").Value = "model_4_8_16"
$ws.Range("B10").Value = 0.6274685010225862
$ws.Range("C10").Value = -0.1292199468941047
$ws.Range("D10").Value = 0.05051557756453018
$ws.Range("E10").Value = -0.042999867852338
$ws.Range("F10").Value = 0.4122826159000397
$ws.Range("G10").Value = 1.472929358482361
$ws.Range("H10").Value = 1.076851844787598
$ws.Range("I10").Value = 1.286540031433105

$ws.Range("A11").Value = "model_4_8_19"
$ws.Range("B11").Value = 0.6281181279095499
$ws.Range("C11").Value = -0.1315074809533789
$ws.Range("D11").Value = 0.05568704152164472
$ws.Range("E11").Value = -0.04204287613716562
$ws.Range("F11").Value = 0.4115636646747589
$ws.Range("G11").Value = 1.475913166999817
$ws.Range("H11").Value = 1.070986747741699
$ws.Range("I11").Value = 1.285359501838684

$ws.Range("A12").Value = "model_4_8_18"
$ws.Range("B12").Value = 0.6291086981772611
$ws.Range("C12").Value = -0.1298785772620183
$ws.Range("D12").Value = 0.06112332484217819
$ws.Range("E12").Value = -0.03877869882098217
$ws.Range("F12").Value = 0.4104673564434052
$ws.Range("G12").Value = 1.473788499832153
$ws.Range("H12").Value = 1.064821124076843
$ws.Range("I12").Value = 1.281333208084106

$ws.Range("A13").Value = "model_4_8_20"
$ws.Range("B13").Value = 0.6293786538106247
$ws.Range("C13").Value = -0.1321652265791273
$ws.Range("D13").Value = 0.064238876050562
$ws.Range("E13").Value = -0.03871083573437906
$ws.Range("F13").Value = 0.4101686179637909
$ws.Range("G13").Value = 1.476771116256714
$ws.Range("H13").Value = 1.061287641525269
$ws.Range("I13").Value = 1.281249284744263

$ws.Range("A14").Value = "model_4_8_14"
$ws.Range("B14").Value = 0.6314278531976811
$ws.Range("C14").Value = -0.1177641745179763
$ws.Range("D14").Value = 0.06472593760201317
$ws.Range("E14").Value = -0.03043764667925686
$ws.Range("F14").Value = 0.4079007506370544
$ws.Range("G14").Value = 1.45798671245575
$ws.Range("H14").Value = 1.06073522567749
$ws.Range("I14").Value = 1.271044611930847

$ws.Range("A15").Value = "model_4_8_17"
$ws.Range("B15").Value = 0.6316389457333789
$ws.Range("C15").Value = -0.1211546863532611
$ws.Range("D15").Value = 0.06766731204835452
$ws.Range("E15").Value = -0.03106335346953859
$ws.Range("F15").Value = 0.4076671600341797
$ws.Range("G15").Value = 1.462409257888794
$ws.Range("H15").Value = 1.057399272918701
$ws.Range("I15").Value = 1.271816253662109

$ws.Range("A16").Value = "model_4_8_8"
$ws.Range("B16").Value = 0.6406883872235041
$ws.Range("C16").Value = -0.07119500867106554
$ws.Range("D16").Value = 0.09183769922039875
$ws.Range("E16").Value = 0.007363762174196764
$ws.Range("F16").Value = 0.3976520895957947
$ws.Range("G16").Value = 1.397242903709412
$ws.Range("H16").Value = 1.029986619949341
$ws.Range("I16").Value = 1.224416375160217

$ws.Range("A17").Value = "model_4_8_7"
$ws.Range("B17").Value = 0.6408904560059603
$ws.Range("C17").Value = -0.06928780763364406
$ws.Range("D17").Value = 0.09197539625815254
$ws.Range("E17").Value = 0.008491222385655472
$ws.Range("F17").Value = 0.3974284529685974
$ws.Range("G17").Value = 1.394755244255066
$ws.Range("H17").Value = 1.029830455780029
$ws.Range("I17").Value = 1.223025798797607

$ws.Range("A18").Value = "model_4_8_13"
$ws.Range("B18").Value = 0.6412870941971612
$ws.Range("C18").Value = -0.0922158452837003
$ws.Range("D18").Value = 0.101767543426714
$ws.Range("E18").Value = -0.0001081012309180895
$ws.Range("F18").Value = 0.396989494562149
$ws.Range("G18").Value = 1.424662113189697
$ws.Range("H18").Value = 1.018724679946899
$ws.Range("I18").Value = 1.233633041381836

$ws.Range("A19").Value = "model_4_8_12"
$ws.Range("B19").Value = 0.6415865160413092
$ws.Range("C19").Value = -0.07749238181008877
$ws.Range("D19").Value = 0.09129404248846296
$ws.Range("E19").Value = 0.003603083181670796
$ws.Range("F19").Value = 0.3966580927371979
$ws.Range("G19").Value = 1.405457019805908
$ws.Range("H19").Value = 1.030603170394897
$ws.Range("I19").Value = 1.229055285453796

$ws.Range("A20").Value = "model_4_8_6"
$ws.Range("B20").Value = 0.6502616257348466
$ws.Range("C20").Value = -0.02434665208738496
$ws.Range("D20").Value = 0.1188551235510156
$ws.Range("E20").Value = 0.04528098975513428
$ws.Range("F20").Value = 0.3870573043823242
$ws.Range("G20").Value = 1.336135029792786
$ws.Range("H20").Value = 0.9993449449539185
$ws.Range("I20").Value = 1.177645683288574

$ws.Range("A21").Value = "model_4_8_5"
$ws.Range("B21").Value = 0.6540721704436856
$ws.Range("C21").Value = 0.04539025176216283
$ws.Range("D21").Value = 0.1059703810479397
$ws.Range("E21").Value = 0.07874698604117536
$ws.Range("F21").Value = 0.3828401565551758
$ws.Range("G21").Value = 1.245171785354614
$ws.Range("H21").Value = 1.013958215713501
$ws.Range("I21").Value = 1.136365175247192

$ws.Range("A22").Value = "model_4_8_4"
$ws.Range("B22").Value = 0.6542515067731459
$ws.Range("C22").Value = 0.1635154596462454
$ws.Range("D22").Value = 0.03309035297881335
$ws.Range("E22").Value = 0.1133433531172102
$ws.Range("F22").Value = 0.3826417326927185
$ws.Range("G22").Value = 1.091091752052307
$ws.Range("H22").Value = 1.096614599227905
$ws.Range("I22").Value = 1.093690633773804

$ws.Range("A23").Value = "model_4_8_9"
$ws.Range("B23").Value = 0.6570777707352462
$ws.Range("C23").Value = -0.02763427200759483
$ws.Range("D23").Value = 0.1412597755639344
$ws.Range("E23").Value = 0.05313462448520401
$ws.Range("F23").Value = 0.3795138895511627
$ws.Range("G23").Value = 1.340423226356506
$ws.Range("H23").Value = 0.9739348888397217
$ws.Range("I23").Value = 1.16795825958252

$ws.Range("A24").Value = "model_4_8_10"
$ws.Range("B24").Value = 0.6639069454294405
$ws.Range("C24").Value = -0.006725342956210634
$ws.Range("D24").Value = 0.1599633217782813
$ws.Range("E24").Value = 0.07293275088286355
$ws.Range("F24").Value = 0.371955931186676
$ws.Range("G24").Value = 1.31315004825592
$ws.Range("H24").Value = 0.9527223110198975
$ws.Range("I24").Value = 1.143537163734436

$ws.Range("A25").Value = "model_4_8_11"
$ws.Range("B25").Value = 0.6645320791525886
$ws.Range("C25").Value = -0.009244097743175361
$ws.Range("D25").Value = 0.1661670159683831
$ws.Range("E25").Value = 0.07420691634656817
$ws.Range("F25").Value = 0.3712641596794128
$ws.Range("G25").Value = 1.31643545627594
$ws.Range("H25").Value = 0.9456865787506104
$ws.Range("I25").Value = 1.141965389251709

$ws.Range("A26").Value = "model_4_8_3"
$ws.Range("B26").Value = 0.7078237969766707
$ws.Range("C26").Value = 0.4438762907675793
$ws.Range("D26").Value = 0.2977487417073772
$ws.Range("E26").Value = 0.3848116229763152
$ws.Range("F26").Value = 0.3233529627323151
$ws.Range("G26").Value = 0.7253954410552979
$ws.Range("H26").Value = 0.7964538931846619
$ws.Range("I26").Value = 0.7588346600532532
